# AFA 2020.xlsx - Finals MI vs DC results entry
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1. Sheet1 - fill in the Eliminator / Qualifier 2 / Finals match rows
# ---------------------------------------------------------------------------

# Row 67 = Eliminator (SRH vs RCB) -- Match text already present (C67)
$ws1.Range("E67").Value = 20
$ws1.Range("H67").Value = 0
$ws1.Range("K67").Value = 40
$ws1.Range("N67").Value = 80
$ws1.Range("Q67").Value = 100
$ws1.Range("T67").Value = 60

# Row 68 = Qualifier 2 (DC vs SRH)
$ws1.Range("C68").Value = "DC vs SRH"
$ws1.Range("E68").Value = 60
$ws1.Range("H68").Value = 80
$ws1.Range("K68").Value = 20
$ws1.Range("N68").Value = 100
$ws1.Range("Q68").Value = 40
$ws1.Range("T68").Value = 0

# Row 69 = Finals (MI vs DC)
$ws1.Range("C69").Value = "MI vs DC"
$ws1.Range("E69").Value = 100
$ws1.Range("H69").Value = 40
$ws1.Range("K69").Value = 60
$ws1.Range("N69").Value = 80
$ws1.Range("Q69").Value = 0
$ws1.Range("T69").Value = 20

# ---------------------------------------------------------------------------
# 2. Remove the blank spacer row (old row 70) - everything below shifts up
# ---------------------------------------------------------------------------
$ws1.Rows(70).Delete()

# ---------------------------------------------------------------------------
# 3. Scorecard (Qualifier1/Eliminator/Qualifier2/Finals win counts) - rows
#    80-85 after the shift (was 81-86)
# ---------------------------------------------------------------------------
# Anantha
$ws1.Range("D80").Value = 3
$ws1.Range("E80").Value = 5
$ws1.Range("F80").Value = 3
# Jayanth
$ws1.Range("D81").Value = 3
$ws1.Range("E81").Value = 10
$ws1.Range("F81").Value = 3
# Justin
$ws1.Range("C82").Value = 10
$ws1.Range("D82").Value = 0
$ws1.Range("E82").Value = 5
$ws1.Range("F82").Value = 3
# Rapaka
$ws1.Range("C83").Value = 7
$ws1.Range("D83").Value = 0
$ws1.Range("E83").Value = 3
$ws1.Range("F83").Value = 0
# Sushma
$ws1.Range("C84").Value = 3
$ws1.Range("D84").Value = 0
$ws1.Range("E84").Value = 0
$ws1.Range("F84").Value = 3
# Sampath M
$ws1.Range("C85").Value = 3
$ws1.Range("D85").Value = 3
$ws1.Range("E85").Value = 7
$ws1.Range("F85").Value = 3

# "Congrats" tag next to the 3 people with a positive final payout
$ws1.Range("S81").Value = "Congrats"
$ws1.Range("S82").Value = "Congrats"
$ws1.Range("S83").Value = "Congrats"
